$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Prova1")
$ws1.Range("A3").Value = "Rosa"
$ws1.Range("B3").Value = "HOJE"
$ws1.Range("H3").Value = "00:00:05"

$ws2 = $wb.Worksheets.Item("Prova2")
$ws2.Range("A3").Value = "Rosa"
$ws2.Range("B3").Value = "HOJE"
$ws2.Range("H3").Value = "00:00:10"
